$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp update (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 12:22"

# --- Simple numeric updates (city stays on the same row) ---
# Row 4: Madrid
$ws.Range("B4").Value = 47146
$ws.Range("C4").Value = 25385
$ws.Range("D4").Value = 15338
$ws.Range("E4").Value = 6423

# Row 10: Navarra
$ws.Range("B10").Value = 4092
$ws.Range("C10").Value = 695
$ws.Range("D10").Value = 3158
$ws.Range("E10").Value = 239

# Row 12: La Rioja
$ws.Range("B12").Value = 3358
$ws.Range("C12").Value = 1458
$ws.Range("D12").Value = 1677
$ws.Range("E12").Value = 223

# Row 47: Huesca
$ws.Range("B47").Value = 505
$ws.Range("C47").Value = 106
$ws.Range("D47").Value = 331
$ws.Range("E47").Value = 68

# --- Rows that swap ranking order (city names trade rows) plus updated counts ---

# Rows 13/14: Alacant/Alicante and Zaragoza swap places; Zaragoza (now row13) gets new counts,
# Alacant/Alicante (now row14) keeps its prior counts.
$ws.Range("A13").Value = "Zaragoza"
$ws.Range("B13").Value = 3137
$ws.Range("C13").Value = 717
$ws.Range("D13").Value = 2075
$ws.Range("E13").Value = 345

$ws.Range("A14").Value = "Alacant/Alicante"
$ws.Range("B14").Value = 3133
$ws.Range("C14").Value = 990
$ws.Range("D14").Value = 1812
$ws.Range("E14").Value = 331

# Rows 26/27: Granada and Cantabria swap places; Cantabria (now row26) gets new counts,
# Granada (now row27) keeps its prior counts.
$ws.Range("A26").Value = "Cantabria"
$ws.Range("B26").Value = 1777
$ws.Range("C26").Value = 317
$ws.Range("D26").Value = 1343
$ws.Range("E26").Value = 117

$ws.Range("A27").Value = "Granada"
$ws.Range("B27").Value = 1772
$ws.Range("C27").Value = 317
$ws.Range("D27").Value = 1290
$ws.Range("E27").Value = 165

# Rows 48/49: Gran Canaria and Teruel swap places; Teruel (now row48) gets new counts,
# Gran Canaria (now row49) keeps its prior counts.
$ws.Range("A48").Value = "Teruel"
$ws.Range("B48").Value = 496
$ws.Range("C48").Value = 113
$ws.Range("D48").Value = 333
$ws.Range("E48").Value = 50

$ws.Range("A49").Value = "Gran Canaria"
$ws.Range("B49").Value = 481
$ws.Range("C49").Value = 120
$ws.Range("D49").Value = 334
$ws.Range("E49").Value = 27
